$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.002.32'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.908.88'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7907'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.72'
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3167'
$ws.Range("E8").Value = '  +3.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.30'
$ws.Range("E9").Value = '  +3.77%  '
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07999'
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.904.57'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7423'
$ws.Range("E13").Value = '  -1.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.192'
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.98'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").Value = '30.008.31'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.91'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.867'
$ws.Range("E18").Value = '  -4.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.26'
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007731'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").Value = '2.143.08'
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.829'
$ws.Range("E24").Value = '  -3.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '168.15'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.221'
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1391'
$ws.Range("E27").Value = '  +10.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.88'
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.035'
$ws.Range("E29").Value = '  -0.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.365'
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.518'
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.310'
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.082'
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05515'
$ws.Range("E34").Value = '  +2.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.254'
$ws.Range("E35").Value = '  -2.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7332'
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01926'
$ws.Range("E38").Value = '  -0.89%  '
$ws.Range("E39").Value = '  +0.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.137'
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4418'
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.21'
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8375'
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.873'
$ws.Range("E45").Value = '  -3.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.55'
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.537'
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '986.08'
$ws.Range("E48").Value = '  +8.30%  '
$ws.Range("D49").Value = '2.052.08'
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.22'
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.477'
$ws.Range("E51").Value = '  +0.28%  '
